$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 97
$prev = 96

# Copy formatting (styles) from the previous row's A/E cells so the new
# row re-uses the existing cellXfs entries (bold/border/centered for A,
# the custom date-time number format for E) instead of synthesising new
# style records.
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$ws.Cells.Item($prev, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item($prev, 5).Copy() | Out-Null
$ws.Cells.Item($row, 5).PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# Values for the new match row (index 96 -> Resistencia vs Sp. Luqueno).
$ws.Cells.Item($row, 1).Value = 96
$ws.Cells.Item($row, 2).Value = "paraguay"
$ws.Cells.Item($row, 3).Value = "primera-division"
$ws.Cells.Item($row, 4).Value = "'2023"
$ws.Cells.Item($row, 5).Value = 45225.02083333334
$ws.Cells.Item($row, 6).Value = "Resistencia"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Sp. Luqueno"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 3.27
$ws.Cells.Item($row, 11).Value = "22/10/2023 00:42"
$ws.Cells.Item($row, 12).Value = 3.97
$ws.Cells.Item($row, 13).Value = "26/10/2023 00:00"
$ws.Cells.Item($row, 14).Value = 3.46
$ws.Cells.Item($row, 15).Value = "22/10/2023 00:42"
$ws.Cells.Item($row, 16).Value = 3.51
$ws.Cells.Item($row, 17).Value = "26/10/2023 00:00"
$ws.Cells.Item($row, 18).Value = 2.14
$ws.Cells.Item($row, 19).Value = "22/10/2023 00:42"
$ws.Cells.Item($row, 20).Value = 2.01
$ws.Cells.Item($row, 21).Value = "26/10/2023 00:00"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/resistencia-sp-luqueno/jkZNTDW6/"
